# Applies the data refresh to the "Applied Torque" summary sheet:
# new test readings, a new calibration serial/model/unit number, and a
# new customer (MASTEC CANADA INC. replaces TIRE CRAFT), plus clearing
# the address column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep every value written below as literal text (matching the
# workbook's existing inlineStr/text-typed cells) instead of letting
# Excel auto-convert numeric-looking strings into real numbers.
$ws.Range("A2:P4").NumberFormat = "@"

$rows = @(2, 3, 4)

$testValues = @{
    2 = @("221.6", "223.8", "225.7", "228.5", "227.7")
    3 = @("147.1", "149.1", "149.2", "148.6", "148.5")
    4 = @("67.9", "68.0", "67.3", "67.4", "70.0")
}

foreach ($r in $rows) {
    $vals = $testValues[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]   # C - Test 1
    $ws.Cells.Item($r, 4).Value = $vals[1]   # D - Test 2
    $ws.Cells.Item($r, 5).Value = $vals[2]   # E - Test 3
    $ws.Cells.Item($r, 6).Value = $vals[3]   # F - Test 4
    $ws.Cells.Item($r, 7).Value = $vals[4]   # G - Test 5

    $ws.Cells.Item($r, 9).Value = "2286916218"          # I - Serial Number
    $ws.Cells.Item($r, 10).Value = "718976"             # J - Model
    $ws.Cells.Item($r, 13).Value = "TW-78"              # M - Unit Number
    $ws.Cells.Item($r, 14).Value = "MASTEC CANADA INC." # N - Customer/Company
    $ws.Cells.Item($r, 15).Value = "403-852-5420"       # O - Phone Number
    $ws.Range("P" + $r).ClearContents()                 # P - Address (cleared)
}
